$wb = $excel.ActiveWorkbook

# "Overview" sheet - Latest HO Xliff Generate Date for a28ccbb7-50ac-4278-80dd-549770ed3b63.md
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-16 15:06:42"

# "zh-cn" sheet - Correspond Handoff / Handback Datetime for a28ccbb7... row
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-16 15:06:36"
$wsZhCn.Range("K2").Value = "2016-08-16 15:06:56"

# "de-de" sheet - Correspond Handoff Datetime for a28ccbb7... row
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-16 15:07:13"
